$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new sign-up entries (OSI Svomming accounting paperwork members)
$ws.Range("B29").Value = "Manu"
$ws.Range("B30").Value = "Hannah"

# Replace "Anne Martina" with "Anne Kraus" everywhere she appears (passenger slots)
$ws.Range("D8").Value = "Anne Kraus"
$ws.Range("B21").Value = "Anne Kraus"
$ws.Range("C22").Value = "Anne Kraus"

$ws.Range("B31").Value = "Louis"
$ws.Range("B32").Value = "Bjørn"

# Move active selection to D13, as last seen in the file
$ws.Range("D13").Select()
